# data: testing_data.xlsx: rename Hamburguer -> Egg
$wb = $excel.ActiveWorkbook

# Rename the "Hamburger" food item to "Egg" on the "foods" sheet (B2).
$wsFoods = $wb.Worksheets.Item("foods")
$wsFoods.Range("B2").Value = "Egg"

# The "foods" sheet's selection moves to C4 and it is no longer the active tab.
$wsFoods.Range("C4").Select()

# The "foods_nutrients" sheet becomes the active tab with selection F9.
$wsFoodsNutrients = $wb.Worksheets.Item("foods_nutrients")
$wsFoodsNutrients.Activate()
$wsFoodsNutrients.Range("F9").Select()
